$d = $word.ActiveDocument

# --- Edit 1: merge the three runs "<id>", "p141r_1", "</id>" into a single
# run carrying the text "<id>p141r_1</id>" (keeping the first run's
# Courier-New / 7f6000 / sz18 formatting, which Word naturally applies to
# merged text). We go through a throw-away marker first so the Find/Replace
# range actually differs from the original text, forcing Word to coalesce
# the three runs into one when the new text is written back in a second
# pass.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("<id>p141r_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Text = "IRON_TMP_MERGE_MARKER_1"
}

$rng1b = $d.Content
$found1b = $rng1b.Find.Execute("IRON_TMP_MERGE_MARKER_1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1b) {
    $rng1b.Text = "<id>p141r_1</id>"
}

# --- Edit 2: remove the comment (and its commentRangeStart/commentRangeEnd/
# commentReference markers) that wraps the single letter "C" in the body
# text, together with the comment content itself. Walk backwards since
# deleting shifts indices.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}
